# Update column AL (18/04/2020) values for the relevant state rows on the
# "20th April 1st update" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    3  = 2    # Andaman and Nicobar Islands
    7  = 1    # Bihar
    8  = 2    # Chandigarh
    10 = 186  # Delhi
    12 = 277  # Gujarat (was 173)
    13 = 9    # Haryana
    14 = 1    # Himachal Pradesh
    15 = 13   # Jammu and Kashmir
    16 = 1    # Jharkhand
    17 = 25   # Karnataka
    18 = 4    # Kerala
    20 = 92   # Madhya Pradesh (was 45)
    21 = 328  # Maharashtra
    26 = 1    # Odisha
    28 = 23   # Punjab
    29 = 122  # Rajasthan (was 41)
    30 = 49   # Tamil Nadu
    31 = 43   # Telangana
    33 = 125  # Uttar Pradesh
    34 = 2    # Uttarakhand
}

foreach ($row in $updates.Keys) {
    $ws.Range("AL$row").Value = $updates[$row]
}
